$d = $word.ActiveDocument

# 1. Merge split title run "H.S.P(Healthy Soft" + "ware Products)" into one run.
#    This Find/Replace also removes the now-orphaned _GoBack bookmark that sat
#    between the two runs.
$d.Content.Find.Execute("H.S.P(Healthy Soft" + "ware Products)", $false, $false, $false, $false, $false, $true, 1, $false, "H.S.P(Healthy Software Products)", 2) | Out-Null

# 2. Merge the NIH paragraph's two runs into one (identical text, just un-split).
$d.Content.Find.Execute("a la transformación y se realizarán", $false, $false, $false, $false, $false, $true, 1, $false, "a la transformación y se realizarán", 2) | Out-Null

# 5. Merge "Ofrecer una herramienta..." runs.
$d.Content.Find.Execute("imágenes/videos/tutoriales", $false, $false, $false, $false, $false, $true, 1, $false, "imágenes/videos/tutoriales", 2) | Out-Null

# 6. Merge "Realizar módulos..." runs.
$d.Content.Find.Execute("banco de sangre, inscripción", $false, $false, $false, $false, $false, $true, 1, $false, "banco de sangre, inscripción", 2) | Out-Null

# 7. Merge "Desarrollar un repositorio..." runs.
$d.Content.Find.Execute("repositorio contará con varias fases", $false, $false, $false, $false, $false, $true, 1, $false, "repositorio contará con varias fases", 2) | Out-Null

# 8. Merge "El proyecto se enfoca..." runs.
$d.Content.Find.Execute("calidad que optimice los procesos", $false, $false, $false, $false, $false, $true, 1, $false, "calidad que optimice los procesos", 2) | Out-Null

# 9. Merge "La dirección del proyecto..." (three runs -> one).
$d.Content.Find.Execute("el líder técnico (Gabriel", $false, $false, $false, $false, $false, $true, 1, $false, "el líder técnico (Gabriel", 2) | Out-Null
$d.Content.Find.Execute("Gabriel Santiago Álvarez Amaya,  Camilo", $false, $false, $false, $false, $false, $true, 1, $false, "Gabriel Santiago Álvarez Amaya,  Camilo", 2) | Out-Null

# 10. Merge "Todas las solicitudes..." runs.
$d.Content.Find.Execute("en tiempo y costo del proyecto", $false, $false, $false, $false, $false, $true, 1, $false, "en tiempo y costo del proyecto", 2) | Out-Null

# 11. Merge "Administrador DB..." runs.
$d.Content.Find.Execute("Jair Dario Muñoz Aguilar)", $false, $false, $false, $false, $false, $true, 1, $false, "Jair Dario Muñoz Aguilar)", 2) | Out-Null

Write-Output "done"
